# Project "Sample Project" is saved.
# The rule table's last row label (B11, was "R40") is corrected to "1".
# A leading apostrophe forces Excel to store the numeric-looking entry as
# text (shared string) rather than auto-converting it to the number 1,
# matching the original author's intent of a literal text value "1".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Rules")

$ws.Range("B11").Value = "'1"
